$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.543.45"
$ws.Range("E2").Value = "  -2.55%  "

$ws.Range("D3").Value = "2.368.68"
$ws.Range("E3").Value = "  -4.14%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.29%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -3.78%  "

$ws.Range("E10").Value = "  -3.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.23"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -8.41%  "

$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("D13").Value = "2.730.41"
$ws.Range("E13").Value = "  -4.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.53"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.06"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.58%  "

$ws.Range("D16").Value = "2.349.94"
$ws.Range("E16").Value = "  -5.03%  "

$ws.Range("E17").Value = "  -4.63%  "

$ws.Range("D18").Value = "40.454.43"
$ws.Range("E18").Value = "  -2.63%  "

$ws.Range("D19").Value = "0.0₃0909"
$ws.Range("E19").Value = "  -3.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("D20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.57"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.41%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  -8.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.03%  "

$ws.Range("E28").Value = "  -2.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.25"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.98"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.10%  "

$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0726"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.70%  "

$ws.Range("E35").Value = "  -5.67%  "

$ws.Range("E36").Value = "  -2.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.15"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.49%  "

$ws.Range("E38").Value = "  -3.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.74"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.89%  "

$ws.Range("E42").Value = "  -5.82%  "

$ws.Range("D43").Value = "1.956.42"
$ws.Range("E43").Value = "  -1.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.74"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.69%  "

$ws.Range("E47").Value = "  -9.14%  "

$ws.Range("D48").Value = "2.599.25"
$ws.Range("E48").Value = "  -3.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.81"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.90"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.40%  "
